$d = $word.ActiveDocument

# Collapse to the very end of the document, right after "Thank you for playing!"
$r = $d.Content
$r.Collapse(0)          # wdCollapseEnd
$r.InsertParagraphAfter()   # adds the blank paragraph

# Re-fetch the end-of-document range (prior range objects don't auto-track
# the new content) and add the "Synopsys" heading paragraph.
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Content
$r.Collapse(0)
$r.InsertAfter("Synopsys")

# Finally, add the synopsis paragraph itself.
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Content
$r.Collapse(0)
$r.InsertAfter("A group of space frogs have been exiled from their planet, and are seeking refuge. Lend a hand for these poor hapless creatures, and they shall be forever grateful. Make use of your knowledge about climate and weather to find them a suitable home across Earth!")
